# "break out stock.yaml completed" -- append the newly-scraped BSE/CONCOR/ZEEL
# rows (24/06/2024 06:44:46 run) to the "10per change" sheet, and fix up the
# bsecode column for the previous (05:45:27) run so it is stored as a real
# number instead of text, matching the upstream screener export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# --- bsecode (column E) for rows 32-34 was imported as text; normalize to numeric ---
$ws.Range("E32").Value = 20
$ws.Range("E33").Value = 531344
$ws.Range("E34").Value = 505537

# --- newly broken-out rows from the 24/06/2024 06:44:46 screener run ---
$ws.Range("A35").Value = "24/06/2024 06:44:46"
$ws.Range("B35").Value = 1
$ws.Range("C35").Value = "BSE"
$ws.Range("D35").Value = "BSE (Bombay stock exchange)"
$ws.Range("E35").Value = "'20"
$ws.Range("F35").Value = -2.35
$ws.Range("G35").Value = 2500
$ws.Range("H35").Value = 572158

$ws.Range("A36").Value = "24/06/2024 06:44:46"
$ws.Range("B36").Value = 2
$ws.Range("C36").Value = "CONCOR"
$ws.Range("D36").Value = "Container Corporation Of India Limited"
$ws.Range("E36").Value = "'531344"
$ws.Range("F36").Value = -3.04
$ws.Range("G36").Value = 1057.7
$ws.Range("H36").Value = 2695405

$ws.Range("A37").Value = "24/06/2024 06:44:46"
$ws.Range("B37").Value = 3
$ws.Range("C37").Value = "ZEEL"
$ws.Range("D37").Value = "Zee Entertainment Enterprises Limited"
$ws.Range("E37").Value = "'505537"
$ws.Range("F37").Value = -1.56
$ws.Range("G37").Value = 151.84
$ws.Range("H37").Value = 6264136
